$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dc.contributor.author")

# Allow referencing objects by their business identifiers: switch the
# authority-key separator from "::" to "$$" for the affected metadata values.
$ws.Range("B2").Value = "Author1`$`$authority1`$`$xxx"
$ws.Range("C3").Value = "OrgUnit2`$`$authority2`$`$400"

# The edits were made on the "dc.contributor.author" sheet, so it ends up
# being the active tab, with the cursor left on the cell right after the
# last edited one.
$ws.Activate()
$ws.Range("C4").Select()
